$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains one new data row: a new weekly price record is inserted
# right after row 3 (pushing the former rows 4-6 down to rows 5-7 unchanged).
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new record's values.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Femacal de La Calera"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44438
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107002
$ws.Range("J4").Value = "Chirimoya"
$ws.Range("K4").Value = "Cultivar IV Región"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 85
$ws.Range("N4").Value = 27000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 28588
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("R4").Value = "Provincia del Elquí"
$ws.Range("S4").Value = 2859
$ws.Range("T4").Value = 10
